# Sync attendance_reports: reorder "Recorded By" (column G) contributor
# lists on the "Session Analysis Results" sheet so entries that start
# with "System" have that token moved to the end (list order reversed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System,*") {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
        $changed++
    }
}

Write-Output ("Updated " + $changed + " 'Recorded By' cells")
